$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202, shifting existing rows 202:227 down to 203:228
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new record
$ws.Range("A202").Value = 3
$ws.Range("B202").Value = "Femacal de La Calera"
$ws.Range("C202").Value = "Coquimbo"
$ws.Range("D202").Value = 44491
$ws.Range("E202").Value = 5
$ws.Range("F202").Value = 100112040
$ws.Range("G202").Value = "Cilantro"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 120
$ws.Range("K202").Value = 3000
$ws.Range("L202").Value = 3000
$ws.Range("M202").Value = 3000
$ws.Range("N202").Value = '$/docena de atados (3 kilos)'
$ws.Range("O202").Value = "Provincia de Quillota"
$ws.Range("P202").Value = 1000
$ws.Range("Q202").Value = 3
$ws.Range("R202").Value = "Hortaliza"
